# Update the division-answer table: each data cell's text is replaced
# with a new "dividend÷divisor=quotient, remainder" string.
# Cells are addressed by (row, column) in the single table on the page
# rather than by a global text search/replace, because several of the
# old values (e.g. "841÷5=168, 1") appear more than once but map to
# different new values depending on which cell they are in.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($row, $col, $newText) {
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $newText
}

Set-CellText 1 1 "376÷9=41, 7"
Set-CellText 1 2 "546÷4=136, 2"
Set-CellText 1 3 "390÷8=48, 6"
Set-CellText 1 4 "516÷3=172, 0"
Set-CellText 1 5 "793÷7=113, 2"

Set-CellText 5 1 "834÷4=208, 2"
Set-CellText 5 2 "608÷3=202, 2"
Set-CellText 5 3 "427÷3=142, 1"
Set-CellText 5 4 "980÷5=196, 0"
Set-CellText 5 5 "136÷6=22, 4"

Set-CellText 9 1 "180÷3=60, 0"
Set-CellText 9 2 "901÷3=300, 1"
Set-CellText 9 3 "895÷7=127, 6"
Set-CellText 9 4 "168÷7=24, 0"
Set-CellText 9 5 "688÷9=76, 4"

Set-CellText 13 1 "979÷6=163, 1"
Set-CellText 13 2 "382÷5=76, 2"
Set-CellText 13 3 "546÷8=68, 2"
Set-CellText 13 4 "903÷8=112, 7"
Set-CellText 13 5 "200÷8=25, 0"

Set-CellText 17 1 "106÷8=13, 2"
Set-CellText 17 2 "662÷8=82, 6"
Set-CellText 17 3 "491÷2=245, 1"
Set-CellText 17 4 "489÷2=244, 1"
Set-CellText 17 5 "328÷2=164, 0"
